$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.730.45"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.154.13"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.74"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.72"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  +14.89%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.432"
$ws.Range("E10").Value = "  +5.39%  "
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.698.91"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.02"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +5.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.757.91"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  +4.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.153.13"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.08"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.17"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.55"
$ws.Range("E21").Value = "  +5.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.81"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.994"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.522"
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.25"
$ws.Range("E28").Value = "  +12.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0866"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.13"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.31"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("E37").Value = "  +5.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.25"
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.639.61"
$ws.Range("E40").Value = "  +9.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0686"
$ws.Range("E41").Value = "  +2.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.25"
$ws.Range("E42").Value = "  +6.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.85"
$ws.Range("E43").Value = "  +3.28%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.195.35"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("E48").Value = "  +13.60%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.23"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.985"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.42"
$ws.Range("E51").Value = "  +2.68%  "
